$wb = $excel.ActiveWorkbook

# Sheet references
$wsSoftware   = $wb.Worksheets.Item("Software")
$wsExtensions = $wb.Worksheets.Item("Extensions(Ctrl+Shift+X)")

# Add the new shared string / cell value: CloudinaryDotNet in B5 of the Extensions sheet
$wsExtensions.Range("B5").Value = "CloudinaryDotNet"

# Activate the Extensions sheet and select cell B6, matching the new active view state
$wsExtensions.Activate()
$wsExtensions.Range("B6").Select()
